# Natmi following Dr Hou advice
# Update the LR-pairs sheet for Tnfsf11-Tnfrsf11b: add a "Receptor symbol" column
# (previously missing) and expand the data from a single Sending/Target cluster pair
# (FAPs -> FAPs) to the full 2x3 combination of Sending clusters (FAPs, sCs) against
# Target clusters (ECs, FAPs, sCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Tnfsf11"
$ws.Cells.Item(2,3).Value = "Tnfrsf11b"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 1.604474
$ws.Cells.Item(2,8).Value = 4.813422
$ws.Cells.Item(2,9).Value = 0.9077880937234966
$ws.Cells.Item(2,10).Value = 0.9077880937234966
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.1085253333333333
$ws.Cells.Item(2,14).Value = 0.325576
$ws.Cells.Item(2,15).Value = 0.04110788357328589
$ws.Cells.Item(2,16).Value = 0.04110788357328589
$ws.Cells.Item(2,17).Value = 0.1741260756746666
$ws.Cells.Item(2,18).Value = 1.567134681072
$ws.Cells.Item(2,19).Value = 0.03731724726600064
$ws.Cells.Item(2,20).Value = 0.03731724726600064
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Tnfsf11"
$ws.Cells.Item(3,3).Value = "Tnfrsf11b"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 1.604474
$ws.Cells.Item(3,8).Value = 4.813422
$ws.Cells.Item(3,9).Value = 0.9077880937234966
$ws.Cells.Item(3,10).Value = 0.9077880937234966
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 2.214957333333333
$ws.Cells.Item(3,14).Value = 6.644871999999999
$ws.Cells.Item(3,15).Value = 0.8389949644181
$ws.Cells.Item(3,16).Value = 0.8389949644181001
$ws.Cells.Item(3,17).Value = 3.553841452442666
$ws.Cells.Item(3,18).Value = 31.984573071984
$ws.Cells.Item(3,19).Value = 0.7616296393927198
$ws.Cells.Item(3,20).Value = 0.7616296393927199
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Tnfsf11"
$ws.Cells.Item(4,3).Value = "Tnfrsf11b"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 1.604474
$ws.Cells.Item(4,8).Value = 4.813422
$ws.Cells.Item(4,9).Value = 0.9077880937234966
$ws.Cells.Item(4,10).Value = 0.9077880937234966
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.31653
$ws.Cells.Item(4,14).Value = 0.9495899999999999
$ws.Cells.Item(4,15).Value = 0.1198971520086141
$ws.Cells.Item(4,16).Value = 0.1198971520086141
$ws.Cells.Item(4,17).Value = 0.5078641552199999
$ws.Cells.Item(4,18).Value = 4.57077739698
$ws.Cells.Item(4,19).Value = 0.1088412070647761
$ws.Cells.Item(4,20).Value = 0.1088412070647761
$ws.Cells.Item(5,1).Value = "sCs"
$ws.Cells.Item(5,2).Value = "Tnfsf11"
$ws.Cells.Item(5,3).Value = "Tnfrsf11b"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.1629803333333333
$ws.Cells.Item(5,8).Value = 0.488941
$ws.Cells.Item(5,9).Value = 0.09221190627650352
$ws.Cells.Item(5,10).Value = 0.09221190627650352
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.1085253333333333
$ws.Cells.Item(5,14).Value = 0.325576
$ws.Cells.Item(5,15).Value = 0.04110788357328589
$ws.Cells.Item(5,16).Value = 0.04110788357328589
$ws.Cells.Item(5,17).Value = 0.01768749500177778
$ws.Cells.Item(5,18).Value = 0.159187455016
$ws.Cells.Item(5,19).Value = 0.003790636307285257
$ws.Cells.Item(5,20).Value = 0.003790636307285257
$ws.Cells.Item(6,1).Value = "sCs"
$ws.Cells.Item(6,2).Value = "Tnfsf11"
$ws.Cells.Item(6,3).Value = "Tnfrsf11b"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = 0.3333333333333333
$ws.Cells.Item(6,7).Value = 0.1629803333333333
$ws.Cells.Item(6,8).Value = 0.488941
$ws.Cells.Item(6,9).Value = 0.09221190627650352
$ws.Cells.Item(6,10).Value = 0.09221190627650352
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 2.214957333333333
$ws.Cells.Item(6,14).Value = 6.644871999999999
$ws.Cells.Item(6,15).Value = 0.8389949644181
$ws.Cells.Item(6,16).Value = 0.8389949644181001
$ws.Cells.Item(6,17).Value = 0.3609944845057778
$ws.Cells.Item(6,18).Value = 3.248950360552
$ws.Cells.Item(6,19).Value = 0.07736532502538024
$ws.Cells.Item(6,20).Value = 0.07736532502538025
$ws.Cells.Item(7,1).Value = "sCs"
$ws.Cells.Item(7,2).Value = "Tnfsf11"
$ws.Cells.Item(7,3).Value = "Tnfrsf11b"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = 0.3333333333333333
$ws.Cells.Item(7,7).Value = 0.1629803333333333
$ws.Cells.Item(7,8).Value = 0.488941
$ws.Cells.Item(7,9).Value = 0.09221190627650352
$ws.Cells.Item(7,10).Value = 0.09221190627650352
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 0.6666666666666666
$ws.Cells.Item(7,13).Value = 0.31653
$ws.Cells.Item(7,14).Value = 0.9495899999999999
$ws.Cells.Item(7,15).Value = 0.1198971520086141
$ws.Cells.Item(7,16).Value = 0.1198971520086141
$ws.Cells.Item(7,17).Value = 0.05158816491
$ws.Cells.Item(7,18).Value = 0.46429348419
$ws.Cells.Item(7,19).Value = 0.01105594494383802
$ws.Cells.Item(7,20).Value = 0.01105594494383802
